$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per upstream refresh.
# Values in column D that are purely numeric-looking strings are
# written with a leading apostrophe so Excel keeps them as text
# (matching the source data, which stores prices/volumes as text).

$ws.Range("D2").Value = '22.229.43'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").Value = '1.556.29'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").Value = "'288.61"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = "'0.3796"
$ws.Range("E7").Value = '  +1.83%  '
$ws.Range("D8").Value = "'0.3279"
$ws.Range("E8").Value = '  -2.32%  '
$ws.Range("D9").Value = "'44.55"
$ws.Range("E9").Value = '  -7.92%  '
$ws.Range("D10").Value = "'1.136"
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").Value = "'0.07368"
$ws.Range("E11").Value = '  -1.95%  '
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.604.06'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = "'6.756"
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = "'0.00001075"
$ws.Range("E17").Value = '  -4.42%  '
$ws.Range("D18").Value = "'0.06651"
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").Value = "'86.45"
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").Value = "'6.426"
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  -2.38%  '
$ws.Range("D23").Value = "'11.71"
$ws.Range("E23").Value = '  -3.71%  '
$ws.Range("D24").Value = '22.215.24'
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("D25").Value = "'2.286"
$ws.Range("E25").Value = '  -4.91%  '
$ws.Range("D26").Value = "'2.566"
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").Value = "'19.31"
$ws.Range("E28").Value = '  -2.40%  '
$ws.Range("D29").Value = "'4.943"
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '1.770.59'
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = "'122.91"
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").Value = "'1.083"
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("D33").Value = "'5.925"
$ws.Range("E33").Value = '  -4.50%  '
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("D35").Value = "'9.403"
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").Value = "'0.08200"
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("D37").Value = "'0.02353"
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("D38").Value = "'0.06333"
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").Value = "'5.361"
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("D40").Value = "'0.2158"
$ws.Range("E40").Value = '  -5.83%  '
$ws.Range("D41").Value = "'1.239"
$ws.Range("E41").Value = '  -4.89%  '
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").Value = "'0.6072"
$ws.Range("E43").Value = '  -4.49%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = "'13.81"
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("D46").Value = "'3.751"
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = "'0.5897"
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("D48").Value = "'123.11"
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = "'1.974"
$ws.Range("E49").Value = '  -4.64%  '
$ws.Range("D50").Value = "'1.178"
$ws.Range("E50").Value = '  -3.65%  '
$ws.Range("D51").Value = "'0.07067"
$ws.Range("E51").Value = '  -3.01%  '
